# The deck originally carries the "Integral" theme on its (only) slide
# master (ppt/theme/theme1.xml) and the stock "Office Theme" colours on
# the notes master's theme (ppt/theme/theme2.xml). The edit swaps which
# theme is applied where: the slide master (and therefore every slide)
# now uses the default "Office Theme" palette, while "Integral" is kept
# around (as the notes-master theme).
#
# Only the 12-slot colour scheme actually differs between the two themes
# - fontScheme and fmtScheme are byte-identical in both - so re-pointing
# the slide master's theme colour scheme at the Office Theme palette
# reproduces the effective change.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Office Theme colour scheme (RGB packed as r + g*256 + b*65536, i.e.
# the value VBA's RGB() would produce for each hex colour below).
$colors.Item(1).RGB  = 0         # dk1       #000000
$colors.Item(2).RGB  = 16777215  # lt1       #FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2       #44546A
$colors.Item(4).RGB  = 15132391  # lt2       #E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1   #5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2   #ED7D31
$colors.Item(7).RGB  = 10855845  # accent3   #A5A5A5
$colors.Item(8).RGB  = 49407     # accent4   #FFC000
$colors.Item(9).RGB  = 12874308  # accent5   #4472C4
$colors.Item(10).RGB = 4697456   # accent6   #70AD47
$colors.Item(11).RGB = 12673797  # hlink     #0563C1
$colors.Item(12).RGB = 7491477   # folHlink  #954F72

# Best-effort: keep the theme/colour-scheme display names in sync with
# the new palette (harmless no-op on hosts that treat these read-only).
try { $master.Theme.Name = "Office Theme" } catch {}
try { $colors.Name = "Office" } catch {}
